# Generate Report for Handback
#
# Updates the "Latest HO Xliff Generate Date" / handoff / handback
# timestamps for the c8e39951-... file row after a new handback report
# was generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# Row 3 corresponds to c8e39951-57c8-4bf4-8150-3dbc85837702.md
# Column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-21 02:52:14"

# --- zh-cn sheet ------------------------------------------------------
# Row 3 corresponds to c8e39951-57c8-4bf4-8150-3dbc85837702.md
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-21 02:52:10"
$zhcn.Range("K3").Value = "2016-08-21 02:52:27"

# --- de-de sheet ------------------------------------------------------
# Row 3 corresponds to c8e39951-57c8-4bf4-8150-3dbc85837702.md
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-21 02:52:14"
$dede.Range("K3").Value = "2016-08-21 02:52:34"
